$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column F (6), shifting T.C (Azure)/T.C (Desc.)/Error
# columns one to the right, to make room for a new "A.REMARKS" column.
$ws.Columns.Item(6).Insert()

# New header cell for the inserted column.
$hdr = $ws.Cells.Item(1, 6)
$hdr.Value() = "A.REMARKS"

# Give the new header the same look as the neighbouring header cells
# (bold 14pt font, thin left/right border, medium top border, thin bottom border).
$hdr.Font.Bold = $true
$hdr.Font.Size = 14
$hdr.Borders.Item(7).LineStyle = 1
$hdr.Borders.Item(7).Weight = 2
$hdr.Borders.Item(10).LineStyle = 1
$hdr.Borders.Item(10).Weight = 2
$hdr.Borders.Item(8).LineStyle = -4138
$hdr.Borders.Item(8).Weight = -4138
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = 2

# The data rows (2-5) stay empty in the new column, but get a thin
# left/right/top border (open bottom), matching the rest of the table.
$data = $ws.Range("F2:F5")
$data.Borders.Item(7).LineStyle = 1
$data.Borders.Item(7).Weight = 2
$data.Borders.Item(10).LineStyle = 1
$data.Borders.Item(10).Weight = 2
$data.Borders.Item(8).LineStyle = 1
$data.Borders.Item(8).Weight = 2

# Restore the workbook's current selection to match the saved state.
$ws.Range("H21").Select()
